# #249: added test case and updated test file
# Adds a 4-column / 10-row Excel Table ("Table1") to the right of the
# existing picture/shape/chart demo content on Sheet1 (O4:R13), with a
# header row of "Column1".."Column4" and a short comment, matching the
# EPPlus "NvPr" regression-test workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the new table -----------------------------------------
$ws.Range("O4").Value = "Column1"
$ws.Range("P4").Value = "Column2"
$ws.Range("Q4").Value = "Column3"
$ws.Range("R4").Value = "Column4"

# Widen the new columns (O:R) to match the authored column width of 11
# (Excel's ColumnWidth property excludes the ~0.8333 char cell-padding that
# ends up in the stored <col width="..."/>, so back it out here).
$ws.Range("O4:R4").ColumnWidth = 10.166666666666666

# --- Create the table itself -----------------------------------------------
$loRange = $ws.Range("O4:R13")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $loRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.Comment = "Luke, I am your father... seriously..."
$lo.TableStyle = "TableStyleMedium2"

# --- Selection left where the author's session ended up --------------------
$ws.Range("R24").Select()
